$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 700
$ws.Range("J18").Value = 900
$ws.Range("L18").Value = 900
$ws.Range("N18").Value = -1468

$ws.Range("H53").Value = 1498.3914
$ws.Range("I53").Value = 166.66667
$ws.Range("J53").Value = 2354.5
$ws.Range("K53").Value = 166.66667
$ws.Range("L53").Value = 2354.5
$ws.Range("M53").Value = 470.33333
$ws.Range("N53").Value = -3628.5

$ws.Range("H129").Value = 806.7692
$ws.Range("I129").Value = 419.4
$ws.Range("K129").Value = 1258.2
$ws.Range("M129").Value = 3741.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17636.863
$ws.Range("I32").Value = 21159.904
$ws.Range("J32").Value = 4551.2856
$ws.Range("K32").Value = 21159.904
$ws.Range("L32").Value = 4551.2856
$ws.Range("M32").Value = -20872.904
$ws.Range("N32").Value = -5125.2856

$ws.Range("H74").Value = 2217.5
$ws.Range("I74").Value = 2845.3157
$ws.Range("J74").Value = 1299.9231
$ws.Range("K74").Value = 2845.3157
$ws.Range("L74").Value = 1299.9231
$ws.Range("M74").Value = -1971.3157
$ws.Range("N74").Value = -3047.9231

$ws.Range("H77").Value = 2217.5
$ws.Range("I77").Value = 2845.3157
$ws.Range("J77").Value = 1299.9231
$ws.Range("K77").Value = 14226.5785
$ws.Range("L77").Value = 6499.6155
$ws.Range("M77").Value = -9858.5785
$ws.Range("N77").Value = -15235.6155

$ws.Range("H110").Value = 3496.2727
$ws.Range("I110").Value = 2639.8572
$ws.Range("J110").Value = 4995
$ws.Range("K110").Value = 2639.8572
$ws.Range("L110").Value = 4995
$ws.Range("M110").Value = -594.8571999999999
$ws.Range("N110").Value = -9085

$ws.Range("H132").Value = 105462.4
$ws.Range("I132").Value = 6828
$ws.Range("J132").Value = 500000
$ws.Range("K132").Value = 20484
$ws.Range("L132").Value = 1500000
$ws.Range("M132").Value = -17954
$ws.Range("N132").Value = -1505060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 701.6667
$ws.Range("I5").Value = 600
$ws.Range("J5").Value = 752.5
$ws.Range("K5").Value = 600
$ws.Range("L5").Value = 752.5
$ws.Range("M5").Value = -487
$ws.Range("N5").Value = -978.5

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H99").Value = 2495.8
$ws.Range("I99").Value = 3000
$ws.Range("J99").Value = 2369.75
$ws.Range("K99").Value = 3000
$ws.Range("L99").Value = 2369.75
$ws.Range("M99").Value = -1502
$ws.Range("N99").Value = -5365.75

$ws.Range("H134").Value = 37205.17
$ws.Range("I134").Value = 41399.848
$ws.Range("J134").Value = 851.3333
$ws.Range("K134").Value = 124199.544
$ws.Range("L134").Value = 2553.9999
$ws.Range("M134").Value = -121664.544
$ws.Range("N134").Value = -7623.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 168.5
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 333
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 333
$ws.Range("M2").Value = 109
$ws.Range("N2").Value = -559

$ws.Range("H6").Value = 93828340
$ws.Range("I6").Value = 12594000
$ws.Range("K6").Value = 12594000
$ws.Range("M6").Value = -12593887

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("N17").ClearContents()

$ws.Range("H31").Value = 12983.375
$ws.Range("I31").Value = 19421.783
$ws.Range("J31").Value = 4272.5884
$ws.Range("K31").Value = 19421.783
$ws.Range("L31").Value = 4272.5884
$ws.Range("M31").Value = -19126.783
$ws.Range("N31").Value = -4862.5884

$ws.Range("H34").Value = 12983.375
$ws.Range("I34").Value = 19421.783
$ws.Range("J34").Value = 4272.5884
$ws.Range("K34").Value = 19421.783
$ws.Range("L34").Value = 4272.5884
$ws.Range("M34").Value = -19219.783
$ws.Range("N34").Value = -4676.5884

$ws.Range("H41").Value = 16666.666
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 16666.666
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 16666.666
$ws.Range("N41").Value = -17522.666
$ws.Range("M41").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H59").Value = 19333.334
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 21200
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 21200
$ws.Range("M59").Value = -8855
$ws.Range("N59").Value = -23490

$ws.Range("H60").Value = 15200
$ws.Range("J60").Value = 15200
$ws.Range("L60").Value = 15200
$ws.Range("N60").Value = -16222

$ws.Range("H68").Value = 66646.42999999999
$ws.Range("J68").Value = 66646.42999999999
$ws.Range("L68").Value = 66646.42999999999
$ws.Range("N68").Value = -68144.42999999999

$ws.Range("H71").Value = 66646.42999999999
$ws.Range("J71").Value = 66646.42999999999
$ws.Range("L71").Value = 199939.29
$ws.Range("N71").Value = -207427.29

$ws.Range("H74").Value = 33648.375
$ws.Range("J74").Value = 33648.375
$ws.Range("L74").Value = 33648.375
$ws.Range("N74").Value = -35396.375

$ws.Range("H77").Value = 33648.375
$ws.Range("J77").Value = 33648.375
$ws.Range("L77").Value = 100945.125
$ws.Range("N77").Value = -109681.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3484.6191
$ws.Range("J68").Value = 3764.6052
$ws.Range("L68").Value = 11293.8156
$ws.Range("N68").Value = -12915.8156

$ws.Range("H71").Value = 3484.6191
$ws.Range("J71").Value = 3764.6052
$ws.Range("L71").Value = 33881.4468
$ws.Range("N71").Value = -41993.4468

$ws.Range("H81").Value = 4159.7144
$ws.Range("J81").Value = 4159.7144
$ws.Range("L81").Value = 12479.1432
$ws.Range("N81").Value = -14725.1432

$ws.Range("H84").Value = 4159.7144
$ws.Range("J84").Value = 4159.7144
$ws.Range("L84").Value = 37437.4296
$ws.Range("N84").Value = -48669.4296

$ws.Range("H107").Value = 4661.7407
$ws.Range("I107").Value = 25424
$ws.Range("J107").Value = 1050.9131
$ws.Range("K107").Value = 76272
$ws.Range("L107").Value = 3152.7393
$ws.Range("M107").Value = -74352
$ws.Range("N107").Value = -6992.7393

$ws.Range("H131").Value = 139718.86
$ws.Range("J131").Value = 157090.12
$ws.Range("L131").Value = 471270.36
$ws.Range("N131").Value = -481350.36

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 7000
$ws.Range("J52").Value = 7000
$ws.Range("L52").Value = 7000
$ws.Range("N52").Value = -7518

$ws.Range("H102").Value = 5717.5713
$ws.Range("I102").Value = 6253.8335
$ws.Range("K102").Value = 6253.8335
$ws.Range("M102").Value = -4631.8335

$ws.Range("H132").Value = 111626.29
$ws.Range("I132").Value = 104576.9
$ws.Range("J132").Value = 129249.75
$ws.Range("K132").Value = 313730.7
$ws.Range("L132").Value = 387749.25
$ws.Range("M132").Value = -311200.7
$ws.Range("N132").Value = -392809.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H131").Value = 35158
$ws.Range("J131").Value = 35158
$ws.Range("L131").Value = 35158
$ws.Range("N131").Value = -45238

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14833.333
$ws.Range("J54").Value = 14833.333
$ws.Range("L54").Value = 14833.333
$ws.Range("N54").Value = -15873.333

$ws.Range("H136").Value = 1584.2174
$ws.Range("I136").Value = 983.3
$ws.Range("K136").Value = 2949.9
$ws.Range("M136").Value = -399.8999999999996
